# Update countries & provincias Spain
# Applies the refreshed COVID case numbers and re-sorted rank for the
# countries whose totals overtook their neighbour (Canada overtakes
# Belgica, Peru overtakes India) plus routine numeric refreshes for a
# handful of other countries further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 4: Estados Unidos -----------------------------------------
$ws.Range("B4").Value = 976176
$ws.Range("C4").Value = 15525
$ws.Range("E4").Value = 802585
$ws.Range("G4").Value = 702
$ws.Range("H4").Value = 54958

# --- Rows 15-16: Canada overtakes Belgica in rank -------------------
$ws.Range("A15").Value = "Canada"
$ws.Range("B15").Value = 46644
$ws.Range("C15").Value = 1290
$ws.Range("D15").Value = 17239
$ws.Range("E15").Value = 26845
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 95
$ws.Range("H15").Value = 2560

$ws.Range("A16").Value = "Belgica"
$ws.Range("B16").Value = 46134
$ws.Range("C16").Value = 809
$ws.Range("D16").Value = 10785
$ws.Range("E16").Value = 28255
$ws.Range("F16").Value = 891
$ws.Range("G16").Value = 177
$ws.Range("H16").Value = 7094

# --- Rows 19-20: Peru overtakes India in rank -----------------------
$ws.Range("A19").Value = "Peru"
$ws.Range("B19").Value = 27517
$ws.Range("C19").Value = 2186
$ws.Range("D19").Value = 8088
$ws.Range("E19").Value = 18701
$ws.Range("F19").Value = 554
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 728

$ws.Range("A20").Value = "India"
$ws.Range("B20").Value = 26917
$ws.Range("C20").Value = 634
$ws.Range("D20").Value = 5939
$ws.Range("E20").Value = 20152
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 826

# --- Row 53: Sudafrica -----------------------------------------------
$ws.Range("B53").Value = 4546
$ws.Range("C53").Value = 185
$ws.Range("E53").Value = 2986
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 87

# --- Row 54: Egipto ----------------------------------------------------
$ws.Range("B54").Value = 4534
$ws.Range("C54").Value = 215
$ws.Range("D54").Value = 1176
$ws.Range("E54").Value = 3041
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 317

# --- Row 63: Barein ------------------------------------------------
$ws.Range("B63").Value = 2647
$ws.Range("C63").Value = 59
$ws.Range("D63").Value = 1189
$ws.Range("E63").Value = 1450

# --- Row 153: Liechtenstein -----------------------------------------
$ws.Range("B153").Value = 82
$ws.Range("C153").Value = 1
$ws.Range("E153").Value = 26

# --- Row 178: Angola --------------------------------------------------
$ws.Range("B178").Value = 26
$ws.Range("C178").Value = 1
$ws.Range("E178").Value = 18
